$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 314 (high and close columns changed) ---
$ws.Cells.Item(314, 4).Value = 8.0829
$ws.Cells.Item(314, 6).Value = 8.0769

# --- Append new row 315 ---
$ws.Cells.Item(315, 1).Value = 45170.33333333334
$ws.Cells.Item(315, 2).Value = "FX_IDC:USDMOP"
$ws.Cells.Item(315, 3).Value = 8.0771
$ws.Cells.Item(315, 4).Value = 8.0823
$ws.Cells.Item(315, 5).Value = 8.048400000000001
$ws.Cells.Item(315, 6).Value = 8.0654
$ws.Cells.Item(315, 7).Value = 0

# --- Append new row 316 ---
$ws.Cells.Item(316, 1).Value = 45201.375
$ws.Cells.Item(316, 2).Value = "FX_IDC:USDMOP"
$ws.Cells.Item(316, 3).Value = 8.0654
$ws.Cells.Item(316, 4).Value = 8.0684
$ws.Cells.Item(316, 5).Value = 8.0509
$ws.Cells.Item(316, 6).Value = 8.058400000000001
$ws.Cells.Item(316, 7).Value = 0

# --- Append new row 317 ---
$ws.Cells.Item(317, 1).Value = 45231.375
$ws.Cells.Item(317, 2).Value = "FX_IDC:USDMOP"
$ws.Cells.Item(317, 3).Value = 8.0585
$ws.Cells.Item(317, 4).Value = 8.0616
$ws.Cells.Item(317, 5).Value = 8.0381
$ws.Cells.Item(317, 6).Value = 8.043100000000001
$ws.Cells.Item(317, 7).Value = 0

# --- Copy the date-column formatting (style) from row 314's A cell onto the
#     new A-column date cells of rows 315-317, matching the source data's look ---
$ws.Range("A314").Copy()
$ws.Range("A315:A317").PasteSpecial(-4122)
$excel.CutCopyMode = $false
